# "Debug bar is more useful." -- zero out the two machine_readable inputs
# (Total revenues measures / Total program spending measures for the
# 2025-2026 column) and move the on-sheet cursor over to the far side of
# the table so the debug/output columns are in view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("machine_readable")

# Clear the two hard-coded input cells back to 0 -- everything else on
# this sheet (and on "For user (EN)" / "Model" / downstream) is a formula,
# so the whole dependent chain recalculates for free.
$ws.Range("M2").Value = 0
$ws.Range("N3").Value = 0

# Bring the sheet to the front and move the selection/viewport over to the
# debug columns on the right (around Q3) instead of the original A1:L1
# input block.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("Q3").Select()
